# Wellness.xlsx update — add the week of 2025-10-15 entries (rows 478-491)
# following the same layout/format as the existing rows above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Extend formatting from the last existing row (477) down through the new
#    rows (478-491) so the new rows inherit the same cell styles (date format
#    on column A, fonts on B-H, default on I).
# ---------------------------------------------------------------------------
$ws.Range("A477:I477").Copy()
$ws.Range("A478:I491").PasteSpecial(-4122)   # xlPasteFormats

# A handful of rows have no "Localisation douleur" (column G) entry. Those
# cells use a slightly different font style (same as the other blank-G rows
# in the sheet, e.g. row 460) instead of the populated-text style.
foreach ($r in @(479, 482, 486, 490)) {
    $ws.Range("G460").Copy()
    $ws.Range("G$r").PasteSpecial(-4122)     # xlPasteFormats
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Fill in the new rows' data.
#    Columns: A=Date, B=Nom du joueur, C=Volume, D=Intensite, E=Charge,
#             F=Fatigue, G=Localisation douleur, H=Plaisir, I=Charge calc.
# ---------------------------------------------------------------------------
$rows = @(
    @(478, 45945, "Amir Etien",       75, 6, 7, 4, "Pied coup", 0),
    @(479, 45945, "Ilyes Boughanmi",  75, 5, 6, 0, "",          0),
    @(480, 45945, "Romain Thunet",    75, 6, 7, 1, "Jambes",    0),
    @(481, 45945, "Omar Benyounes",   75, 6, 3, 4, "Quadri",   10),
    @(482, 45945, "Naim Ighbane",     75, 5, 3, 0, "",          0),
    @(483, 45945, "Kamal Bafounta",   75, 7, 5, 6, "Genou",     5),
    @(484, 45945, "Maé Clavel",       75, 7, 4, 2, "Ischio",   10),
    @(485, 45945, "Levy Ndoutoume",   75, 7, 7, 1, "Ischio",    4),
    @(486, 45945, "Malik Boussaid",   75, 3, 0, 0, "",         10),
    @(487, 45945, "Emmanuel Valey",   75, 7, 5, 5, "Cheville",  6),
    @(488, 45945, "Karahali Souaré",  75, 5, 4, 6, "Cheville",  1),
    @(489, 45945, "Naim Dhib",        75, 6, 7, 5, "Quadri",    0),
    @(490, 45945, "Sofiane Belle",    75, 5, 3, 0, "",          9),
    @(491, 45945, "Amine Taiar",      75, 4, 3, 3, "Genou",     8)
)

foreach ($row in $rows) {
    $r    = $row[0]
    $date = $row[1]
    $name = $row[2]
    $vol  = $row[3]
    $inte = $row[4]
    $chg  = $row[5]
    $fat  = $row[6]
    $loc  = $row[7]
    $plai = $row[8]

    $ws.Cells.Item($r, 1).Value = $date
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $vol
    $ws.Cells.Item($r, 4).Value = $inte
    $ws.Cells.Item($r, 5).Value = $chg
    $ws.Cells.Item($r, 6).Value = $fat
    if ($loc -ne "") {
        $ws.Cells.Item($r, 7).Value = $loc
    }
    $ws.Cells.Item($r, 8).Value = $plai
    $ws.Cells.Item($r, 9).Formula = "=C" + $r + "*D" + $r
}

# ---------------------------------------------------------------------------
# 3) Leave the view roughly where the author left it (scrolled down, with
#    K486 as the active cell).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("K486").Select()

Write-Host "Added rows 478-491 to Feuil1."
